$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.290.46"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").Value = "1.866.26"

$ws.Range("E4").Value = "  +0.68%  "

$ws.Range("D5").Value = "'239.68"
$ws.Range("E5").Value = "  +3.31%  "

$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  +0.63%  "

$ws.Range("D8").Value = "'42.48"
$ws.Range("E8").Value = "  +6.61%  "

$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").Value = "'0.0696"
$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("D11").Value = "'0.0991"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").Value = "2.134.68"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").Value = "'11.57"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.865.50"
$ws.Range("E14").Value = "  +1.15%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.681"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("E16").Value = "  +1.75%  "

$ws.Range("D17").Value = "35.265.13"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "'70.19"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").Value = "'241.39"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "'12.27"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  +1.40%  "

$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("D24").Value = "'2.26"
$ws.Range("E24").Value = "  -0.60%  "

$ws.Range("D25").Value = "'169.81"

$ws.Range("E26").Value = "  +25.58%  "

$ws.Range("E27").Value = "  +4.24%  "

$ws.Range("D28").Value = "'17.76"
$ws.Range("E28").Value = "  +1.69%  "

$ws.Range("E29").Value = "  +0.74%  "

$ws.Range("E30").Value = "  +1.83%  "

$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("E32").Value = "  +2.14%  "

$ws.Range("E33").Value = "  +27.91%  "

$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").Value = "'2.09"
$ws.Range("E35").Value = "  +8.91%  "

$ws.Range("D36").Value = "'0.821"
$ws.Range("E36").Value = "  +17.92%  "

$ws.Range("E37").Value = "  +6.30%  "

$ws.Range("E38").Value = "  +3.86%  "

$ws.Range("E39").Value = "  +4.30%  "

$ws.Range("D40").Value = "'90.45"
$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("D41").Value = "1.345.70"
$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").Value = "'15.31"
$ws.Range("E42").Value = "  +2.97%  "

$ws.Range("D43").Value = "'0.0603"
$ws.Range("E43").Value = "  +15.53%  "

$ws.Range("E44").Value = "  +2.32%  "

$ws.Range("E45").Value = "  +0.58%  "

$ws.Range("D46").Value = "'12.39"
$ws.Range("E46").Value = "  +46.01%  "

$ws.Range("E47").Value = "  +5.19%  "

$ws.Range("E48").Value = "  -0.79%  "

$ws.Range("D49").Value = "2.051.70"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").Value = "'0.0687"
$ws.Range("E50").Value = "  +3.20%  "

$ws.Range("E51").Value = "  -0.38%  "
